# Update the division problems in the practice-sheet table.
# Each original "old" expression is unique within the document, so a
# simple whole-document Find/ReplaceAll per pair is safe and order-independent.

$d = $word.ActiveDocument

$replacements = @(
    @{old="287÷8="; new="204÷4="},
    @{old="448÷9="; new="589÷2="},
    @{old="674÷5="; new="834÷9="},
    @{old="432÷6="; new="115÷4="},
    @{old="540÷4="; new="672÷7="},
    @{old="249÷8="; new="909÷8="},
    @{old="882÷2="; new="508÷4="},
    @{old="210÷2="; new="662÷3="},
    @{old="745÷6="; new="825÷9="},
    @{old="305÷2="; new="155÷7="},
    @{old="123÷9="; new="569÷8="},
    @{old="776÷3="; new="896÷7="},
    @{old="903÷9="; new="440÷2="},
    @{old="289÷8="; new="213÷9="},
    @{old="670÷9="; new="269÷7="},
    @{old="411÷7="; new="134÷9="},
    @{old="540÷8="; new="666÷9="},
    @{old="413÷4="; new="785÷2="},
    @{old="836÷2="; new="325÷2="},
    @{old="608÷6="; new="866÷8="},
    @{old="155÷2="; new="882÷6="},
    @{old="365÷4="; new="252÷7="},
    @{old="815÷4="; new="846÷9="},
    @{old="606÷5="; new="554÷3="},
    @{old="949÷4="; new="425÷7="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Host "Replaced $($replacements.Count) division expressions"
